$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 769 (shifts existing rows 769..810 down to 770..811)
$ws.Rows.Item(769).Insert()

# Populate the newly inserted row with the new record.
# Column A holds a date-like string that must stay plain text (not be
# auto-converted to a real date by Excel), so force text formatting,
# assign the value, then clear the formatting back to the sheet default
# so the cell ends up unstyled just like its neighbours.
$dateCell = $ws.Cells.Item(769, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026/02/07"
$dateCell.ClearFormats()

$ws.Cells.Item(769, 2).Value = "土"
$ws.Cells.Item(769, 3).Value = 14
$ws.Cells.Item(769, 4).Value = 71
